$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(-1.82958666343958, 0.0353490412064622, -51.757745075844, 0)
    3 = @(0.14798259816539, 0.183539251720133, 0.806272210322832, 0.420085910295573)
    4 = @(0.0067473837056031, 0.0428336452396018, 0.15752532075801, 0.874830853925021)
    5 = @(0.779231371574032, 0.253646555005227, 3.07211494182514, 0.00212547863054266)
    6 = @(0.0508398507941526, 0.0377311858641215, 1.34742255324914, 0.17784418097455)
    7 = @(-0.0494435007468305, 0.0253107681958522, -1.95345713588152, 0.0507654516641872)
    8 = @(0.0269324699716841, 0.0368377989573956, 0.731109641019339, 0.464712185958739)
    9 = @(0.0586438786983845, 0.0359479763407847, 1.6313541030083, 0.102815623459227)
    10 = @(0.0150875548001833, 0.0226845311240007, 0.665103224647229, 0.505984472299242)
    11 = @(0.0847802479361467, 0.0304301698699011, 2.78605897694984, 0.00533531636155563)
    12 = @(-0.125923262455898, 0.0223826877129023, -5.62592232314041, 0.0000000184519370160125)
    13 = @(0.0956058409343752, 0.212600768530514, 0.449696591386749, 0.652929229479633)
    14 = @(-0.0497018743278259, 0.0538796150737326, -0.922461570295379, 0.356287866793132)
    15 = @(-0.604206199242694, 0.260364650911112, -2.32061532596055, 0.0203076142919525)
    16 = @(-0.0740840800291394, 0.0467034208288452, -1.58626667414013, 0.112678828718007)
    17 = @(0.293396047693272, 0.204634407013718, 1.43375716711024, 0.151641572092534)
    18 = @(-0.0598900035319574, 0.0685959376523756, -0.873083823643649, 0.382617390314797)
    19 = @(-0.992759898467075, 0.258839353002161, -3.83542875900631, 0.000125345436177656)
    20 = @(0.163739470069911, 0.0609245637196636, 2.68757722785406, 0.00719724570856946)
    21 = @(0.360790124788675, 0.262839052748353, 1.37266559522303, 0.169856329538468)
    22 = @(0.0379485384244285, 0.0648262248578347, 0.585388683478801, 0.558286396226256)
    23 = @(-0.743704820553769, 0.258975661148473, -2.87171704574739, 0.00408248325473775)
    24 = @(0.00960905261866274, 0.0507156723518069, 0.18946909649558, 0.849725173820309)
    25 = @(0.122598454081093, 0.23319012805072, 0.525744614945394, 0.599065660762274)
    26 = @(-0.0305889369287724, 0.0454605962217586, -0.672867042472525, 0.501031881138791)
    27 = @(-0.570964856066883, 0.26298674926175, -2.17107841999523, 0.0299252449996165)
    28 = @(-0.0562253022150751, 0.0402823342284009, -1.39578063913272, 0.162780562317953)
    29 = @(0.28176112443158, 0.207918650852596, 1.35515079227469, 0.175369518796247)
    30 = @(0.139418126739049, 0.0604692439065695, 2.30560393568774, 0.0211327731768321)
    31 = @(-0.651795641241281, 0.296769383305434, -2.19630352019991, 0.0280702259461882)
    32 = @(0.0594297148253702, 0.0538840639267048, 1.1029182005687, 0.270062689544172)
    33 = @(0.049648714065357, 0.201732150077356, 0.246112055249095, 0.805595491307287)
    34 = @(-0.0821391071271029, 0.0451990432110139, -1.81727535124212, 0.0691749622147985)
    35 = @(-0.709855134600413, 0.265241516797662, -2.67625952064632, 0.00744489718049898)
    36 = @(-0.0107235187642509, 0.0394225543576618, -0.272014813321368, 0.785610629279791)
    37 = @(-0.0375357041319093, 0.0267690236048724, -1.40220669554332, 0.160853532828149)
    38 = @(-0.0395268760045194, 0.0384842876920758, -1.02709127217803, 0.304377488019071)
    39 = @(0.0392668657785028, 0.0384841667654011, 1.02033820864235, 0.307568088261091)
    40 = @(0.00710979050275774, 0.0352615275933371, 0.201630246560877, 0.840205797504941)
    41 = @(-0.0374210608849947, 0.0515041723196673, -0.726563678234378, 0.467493285434933)
    42 = @(-0.071288533576721, 0.053754390372714, -1.3261899741106, 0.184776779368423)
    43 = @(-0.00116699013670652, 0.0264734850079641, -0.0440814700578884, 0.964839463177149)
    44 = @(0.0751498219962275, 0.0389182528739148, 1.93096597218004, 0.0534872606353803)
    45 = @(-0.00242770444859141, 0.0377728865855847, -0.0642710861689318, 0.948754375735482)
    46 = @(-0.247554537569493, 0.262797700684978, -0.941996588722984, 0.346194386412606)
    47 = @(0.0632363052756441, 0.0569360629680677, 1.11065468841971, 0.266717013061332)
    48 = @(0.526569803522911, 0.270175236054214, 1.94899359102345, 0.0512961900983394)
    49 = @(0.061050648551132, 0.0496247501161957, 1.23024596412441, 0.218605012392274)
    50 = @(-0.357772551510803, 0.252519754763098, -1.41681014955225, 0.156538441476267)
    51 = @(0.126941646339557, 0.0723176057460313, 1.75533530224102, 0.0792019714397)
    52 = @(0.667666066568406, 0.268431257417184, 2.48728882393434, 0.0128720848823632)
    53 = @(-0.0248552714642866, 0.0642012194884891, -0.387146407222732, 0.698647823164355)
    54 = @(-0.0616560267215445, 0.317845368590554, -0.193981202227204, 0.846190600780834)
    55 = @(-0.0206409886168359, 0.0691285490432694, -0.298588483376328, 0.765254053740559)
    56 = @(0.348238525950435, 0.268732084618706, 1.29585764366223, 0.195024534205964)
    57 = @(-0.0647317093407376, 0.0542524182998294, -1.19315804473441, 0.232807482468774)
    58 = @(-0.483774404894997, 0.272757696321465, -1.77364162925336, 0.0761224471837897)
    59 = @(0.0392637327469828, 0.0736205244449498, 0.533325903924286, 0.593807999169739)
    60 = @(0.530403921627516, 0.310010041631902, 1.71092497144755, 0.08709496415456)
    61 = @(-0.0731112200284432, 0.0649710207061998, -1.12528969429391, 0.260466295366509)
    62 = @(-0.160471669949967, 0.238126173759586, -0.673893454954602, 0.500379053515439)
    63 = @(0.0188907200759408, 0.096212711546928, 0.196343287411943, 0.844341482252037)
    64 = @(0.721270439114486, 0.305786034222764, 2.3587422523981, 0.0183369869328401)
    65 = @(-0.0785310103475873, 0.0868992853624229, -0.903701451859646, 0.366153729063557)
    67 = @(0.0202637683506314, 0.0949062848033139, 0.21351345058577, 0.830926498946361)
    68 = @(0.757661746231008, 0.306862062141944, 2.46906294294711, 0.0135467388013958)
    69 = @(-0.122963629728301, 0.0749151475177226, -1.64137205628823, 0.100720206813449)
    70 = @(-0.153425861498418, 0.231377140886128, -0.663098614283277, 0.507267398095482)
    71 = @(0.104888209359142, 0.0567515020447037, 1.84820146745226, 0.0645732032267465)
    72 = @(0.624230513373624, 0.272423763993816, 2.29139523007175, 0.0219405652070481)
    73 = @(0.0419978868193393, 0.0487924341680854, 0.860745882746092, 0.389378016593876)
    74 = @(-0.219456684254795, 0.223350644031878, -0.982565710549165, 0.325821227502111)
    75 = @(0.120330121080466, 0.0727783426343069, 1.65337814416983, 0.0982539280899282)
    76 = @(0.769411432197309, 0.270879019441887, 2.84042460646302, 0.00450535212872507)
    77 = @(-0.00635898195139458, 0.0644121611608469, -0.0987233130637434, 0.921357958137152)
    78 = @(-0.220292842169884, 0.278244708117545, -0.791723385002602, 0.428521983786817)
    79 = @(0.0727784012658602, 0.0685470380712443, 1.06172933672521, 0.288358577188707)
    80 = @(0.609516304837646, 0.271205448959503, 2.24743384462257, 0.0246123140669779)
    81 = @(-0.0159833944870081, 0.053282640012634, -0.299973771630277, 0.764197162062067)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 5).Value = $vals[0]
    $ws.Cells.Item($row, 6).Value = $vals[1]
    $ws.Cells.Item($row, 7).Value = $vals[2]
    $ws.Cells.Item($row, 8).Value = $vals[3]
}

$ws.Cells.Item(82, 5).Value = 0.124764479072756

$ws.Rows("83:92").Delete()

Write-Host "edit complete"